# Commit: "added quiz01 and its solutions"
# On the "Load data from a csv into a db" slide, the content placeholder
# (holding the sqlite3/csvsql walkthrough) was nudged from its inherited
# master position to an explicit position/size.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(19)
$sh = $s.Shapes.Item(2)

$sh.Left = 66.0
$sh.Top = 141.2237014874
$sh.Width = 828.0
$sh.Height = 342.6250458701
